$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G9").Value = "SUMPRODUCT"
$ws.Range("H9").Formula = "=SUMPRODUCT(C2:C6,D2:D6)"

$ws.Range("J16").Select()
